$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet/tab to reflect the new "through" date
$ws.Name = "Through 2022-10-21"

# Update the row label for October to reflect the new "through" date
$ws.Range("A11").Value = "October (through 10-21)"

# Update October row (row 11) values for years 2016-2022 (columns C-I)
$ws.Range("C11").Value = 33
$ws.Range("D11").Value = 42
$ws.Range("E11").Value = 50
$ws.Range("F11").Value = 34
$ws.Range("G11").Value = 101
$ws.Range("H11").Value = 130
$ws.Range("I11").Value = 70

# Update Total row (row 12) values for years 2016-2022 (columns C-I)
$ws.Range("C12").Value = 462
$ws.Range("D12").Value = 669
$ws.Range("E12").Value = 598
$ws.Range("F12").Value = 456
$ws.Range("G12").Value = 1002
$ws.Range("H12").Value = 1377
$ws.Range("I12").Value = 1347
